# Auto-upload VRF Excel file
# Adds a new worksheet "qwer" (a duplicate of the "qwe" template sheet)
# at the end of the workbook, with the standard VRF header row.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the very last tab so it lands at the end
# (the default Add() would drop it next to the currently active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "qwer"

$headers = @("Outdoor Model", "Outdoor Quantity", "Outdoor Serial(s)", "Indoor Model", "Indoor Quantity", "Indoor Serial(s)")

for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

$ws.Range("A1").Select() | Out-Null
